$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing rows 2-10 (columns B:G) with new values
$data = @(
    @(2, 0.01778118981931784, 0.8812756614903005, 2.316443170155194, 1.521986586719868, 1.537026199982511, 51),
    @(3, 0.2232283752735154, 0.8890070218121189, 2.324356938520629, 1.524584185448816, 1.5234647779273, 50),
    @(4, 0.07282520553819657, 0.8920146246384424, 2.369955897058822, 1.539466107798032, 1.553678209244919, 49),
    @(5, 0.2280708590847791, 0.9180197143133344, 2.436291467276783, 1.5608624113857, 1.560450066177918, 48),
    @(6, 0.1059889376669026, 0.8775097619739644, 2.394266195833348, 1.547341654526675, 1.560396615019112, 47),
    @(7, 0.2424232887606824, 0.8995544753595199, 2.464361956921573, 1.569828639349395, 1.56813600370666, 46),
    @(8, 0.08228634011495191, 0.8063847497337052, 2.322743373247686, 1.524054911493574, 1.539028286382188, 45),
    @(9, 0.2065927448523342, 0.8408213022771569, 2.429435650446118, 1.558664701097102, 1.562773502624323, 44),
    @(10, 0.1150341483697026, 0.7847972820185366, 2.399366321466991, 1.548988806114167, 1.562992696000677, 43)
)

foreach ($r in $data) {
    $row = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $ws.Cells.Item($row, 6).Value = $r[5]
    $ws.Cells.Item($row, 7).Value = $r[6]
}

# Add new row 11: Q9, copying the header-cell style/formatting from A10
$ws.Range("A10").Copy()
$ws.Range("A11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Cells.Item(11, 1).Value = "Q9"
$ws.Cells.Item(11, 2).Value = 0.124196579601926
$ws.Cells.Item(11, 3).Value = 0.8001464258099134
$ws.Cells.Item(11, 4).Value = 2.371409906460054
$ws.Cells.Item(11, 5).Value = 1.539938280081398
$ws.Cells.Item(11, 6).Value = 1.553527652694611
$ws.Cells.Item(11, 7).Value = 42
